# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values for rows 2-33 with the newly computed strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 3
    3  = 1
    4  = 5
    5  = 2
    6  = 4
    7  = 9
    8  = 8
    9  = 7
    10 = 3
    11 = 7
    12 = 4
    13 = 5
    14 = 6
    15 = 3
    16 = 6
    17 = 6
    18 = 7
    19 = 3
    20 = 3
    21 = 3
    22 = 2
    23 = 3
    24 = 6
    25 = 4
    26 = 6
    27 = 5
    28 = 5
    29 = 5
    30 = 5
    31 = 8
    32 = 1
    33 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
